$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26").Value = 1049.212278876171
$ws.Range("E26").Value = 1716.262336640367
$ws.Range("G26").Value = 418.3870967741935
$ws.Range("H26").Value = 1415.435483870968
$ws.Range("I26").Value = 31
$ws.Range("J26").Value = 2112.168164992067
$ws.Range("K26").Value = 2565.548802598415
$ws.Range("M26").Value = 1281.540983606557
$ws.Range("N26").Value = 3660.180327868853
$ws.Range("O26").Value = 31
$ws.Range("P26").Value = 2110.216475606505
$ws.Range("Q26").Value = 2536.46481911985
$ws.Range("S26").Value = 1537.677685950413
$ws.Range("T26").Value = 3574.628099173554
$ws.Range("U26").Value = 31
$ws.Range("V26").Value = 1826.636559139785
$ws.Range("W26").Value = 2249.948061240582
$ws.Range("Y26").Value = 1325.104761904762
$ws.Range("Z26").Value = 2852.72380952381
$ws.Range("AA26").Value = 31
$ws.Range("AB26").Value = 470.1374877810363
$ws.Range("AC26").Value = 668.4111943299948
$ws.Range("AE26").Value = 313.1575757575758
$ws.Range("AF26").Value = 656.1704545454545
$ws.Range("AG26").Value = 31
$ws.Range("D27").Value = -200.323361082206
$ws.Range("J27").Value = 860.9198836594395
$ws.Range("P27").Value = 873.4358837643297
$ws.Range("V27").Value = 558.1278417818741
$ws.Range("AB27").Value = 44.85631720430109
$ws.Range("D28").Value = 31
$ws.Range("J28").Value = 31
$ws.Range("P28").Value = 31
$ws.Range("V28").Value = 31
$ws.Range("AB28").Value = 31
